$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.510.72'
$ws.Range("E2").Value = '  +6.81%  '

$ws.Range("D3").Value = '1.727.49'
$ws.Range("E3").Value = '  +3.80%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = '332.87'
$ws.Range("E5").Value = '  +0.79%  '

$ws.Range("E6").Value = '  -0.20%  '

$ws.Range("D7").Value = '0.3731'
$ws.Range("E7").Value = '  +2.40%  '

$ws.Range("D8").Value = '48.59'
$ws.Range("E8").Value = '  +2.80%  '

$ws.Range("D9").Value = '0.3377'
$ws.Range("E9").Value = '  +3.84%  '

$ws.Range("D10").Value = '1.185'
$ws.Range("E10").Value = '  +4.24%  '

$ws.Range("D11").Value = '0.07449'
$ws.Range("E11").Value = '  +5.45%  '

$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.07%  '

$ws.Range("D13").Value = '6.408'
$ws.Range("E13").Value = '  +5.68%  '

$ws.Range("D14").Value = '20.15'

$ws.Range("D15").Value = '7.053'
$ws.Range("E15").Value = '  +7.17%  '

$ws.Range("D16").Value = '1.722.35'
$ws.Range("E16").Value = '  +3.39%  '

$ws.Range("D17").Value = '0.00001074'
$ws.Range("E17").Value = '  +2.48%  '

$ws.Range("D18").Value = '0.06661'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("D19").Value = '81.97'
$ws.Range("E19").Value = '  +4.53%  '

$ws.Range("E20").Value = '  -0.09%  '

$ws.Range("D21").Value = '16.56'
$ws.Range("E21").Value = '  +5.04%  '

$ws.Range("D22").Value = '6.161'
$ws.Range("E22").Value = '  +4.09%  '

$ws.Range("D23").Value = '12.78'
$ws.Range("E23").Value = '  +2.05%  '

$ws.Range("D24").Value = '26.517.72'
$ws.Range("E24").Value = '  +6.95%  '

$ws.Range("D25").Value = '2.450'
$ws.Range("E25").Value = '  -0.51%  '

$ws.Range("D26").Value = '1.415'
$ws.Range("E26").Value = '  +21.36%  '

$ws.Range("D27").Value = '2.397'
$ws.Range("E27").Value = '  -1.34%  '

$ws.Range("D28").Value = '151.05'
$ws.Range("E28").Value = '  +1.38%  '

$ws.Range("D29").Value = '19.44'
$ws.Range("E29").Value = '  +4.22%  '

$ws.Range("D30").Value = '1.916.17'
$ws.Range("E30").Value = '  +3.64%  '

$ws.Range("D31").Value = '131.39'
$ws.Range("E31").Value = '  +4.53%  '

$ws.Range("D32").Value = '4.101'
$ws.Range("E32").Value = '  +0.86%  '

$ws.Range("D33").Value = '5.983'
$ws.Range("E33").Value = '  +5.09%  '

$ws.Range("D34").Value = '0.08643'
$ws.Range("E34").Value = '  +1.82%  '

$ws.Range("D35").Value = '1.696'
$ws.Range("E35").Value = '  +3.31%  '

$ws.Range("D36").Value = '12.79'
$ws.Range("E36").Value = '  +5.33%  '

$ws.Range("D37").Value = '5.380'
$ws.Range("E37").Value = '  +4.19%  '

$ws.Range("D38").Value = '0.02333'
$ws.Range("E38").Value = '  +2.30%  '

$ws.Range("D39").Value = '0.06213'
$ws.Range("E39").Value = '  +0.49%  '

$ws.Range("D40").Value = '0.2155'
$ws.Range("E40").Value = '  +3.29%  '

$ws.Range("D41").Value = '8.410'
$ws.Range("E41").Value = '  +2.20%  '

$ws.Range("D42").Value = '1.223'
$ws.Range("E42").Value = '  -1.51%  '

$ws.Range("D43").Value = '0.6212'
$ws.Range("E43").Value = '  +4.90%  '

$ws.Range("D44").Value = '14.20'
$ws.Range("E44").Value = '  +5.90%  '

$ws.Range("D45").Value = '1.000'
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("E46").Value = '  +1.49%  '

$ws.Range("D47").Value = '0.6023'
$ws.Range("E47").Value = '  +6.28%  '

$ws.Range("D48").Value = '128.83'
$ws.Range("E48").Value = '  +2.60%  '

$ws.Range("D49").Value = '2.048'
$ws.Range("E49").Value = '  +5.28%  '

$ws.Range("D50").Value = '0.07186'
$ws.Range("E50").Value = '  +3.13%  '

$ws.Range("D51").Value = '77.01'
$ws.Range("E51").Value = '  +2.57%  '
